# Auto-generated edit script: apply cell-value corrections across 8 leve-profit sheets
# per the commit diff (scheduled runner data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6907.222
$ws.Range("I40").Value = 5148
$ws.Range("K40").Value = 5148
$ws.Range("M40").Value = -4973
$ws.Range("H43").Value = 5894005
$ws.Range("I43").Value = 11114200
$ws.Range("J43").Value = 21286.25
$ws.Range("K43").Value = 11114200
$ws.Range("L43").Value = 21286.25
$ws.Range("M43").Value = -11114131
$ws.Range("N43").Value = -21424.25
$ws.Range("H138").Value = 2867.8
$ws.Range("J138").Value = 2867.8
$ws.Range("L138").Value = 8603.400000000001
$ws.Range("N138").Value = -18883.4
$ws.Range("H139").Value = 32000
$ws.Range("I139").Value = 32000
$ws.Range("K139").Value = 32000
$ws.Range("M139").Value = -26860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2706027.5
$ws.Range("I32").Value = 3058.3076
$ws.Range("K32").Value = 3058.3076
$ws.Range("M32").Value = -2771.3076
$ws.Range("H44").Value = 11772
$ws.Range("J44").Value = 11772
$ws.Range("L44").Value = 11772
$ws.Range("N44").Value = -12748
$ws.Range("H45").Value = 3030.1765
$ws.Range("I45").Value = 2194
$ws.Range("J45").Value = 3970.875
$ws.Range("K45").Value = 2194
$ws.Range("L45").Value = 3970.875
$ws.Range("M45").Value = -1817
$ws.Range("N45").Value = -4724.875
$ws.Range("H61").Value = 4894.3076
$ws.Range("I61").Value = 3866.125
$ws.Range("J61").Value = 6539.4
$ws.Range("K61").Value = 3866.125
$ws.Range("L61").Value = 6539.4
$ws.Range("M61").Value = -3654.125
$ws.Range("N61").Value = -6963.4
$ws.Range("H74").Value = 2193.36
$ws.Range("I74").Value = 1887.3334
$ws.Range("K74").Value = 1887.3334
$ws.Range("M74").Value = -1013.3334
$ws.Range("H77").Value = 2193.36
$ws.Range("I77").Value = 1887.3334
$ws.Range("K77").Value = 9436.666999999999
$ws.Range("M77").Value = -5068.666999999999
$ws.Range("H122").Value = 1711.2413
$ws.Range("I122").Value = 1350.0454
$ws.Range("J122").Value = 2846.4285
$ws.Range("K122").Value = 4050.1362
$ws.Range("L122").Value = 8539.2855
$ws.Range("M122").Value = -1600.1362
$ws.Range("N122").Value = -13439.2855
$ws.Range("H132").Value = 1948.8334
$ws.Range("I132").Value = 1994.6364
$ws.Range("J132").Value = 1445
$ws.Range("K132").Value = 5983.9092
$ws.Range("L132").Value = 4335
$ws.Range("M132").Value = -3453.9092
$ws.Range("N132").Value = -9395
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 4894.3076
$ws.Range("I136").Value = 3866.125
$ws.Range("J136").Value = 6539.4
$ws.Range("K136").Value = 11598.375
$ws.Range("L136").Value = 19618.2
$ws.Range("M136").Value = -9048.375
$ws.Range("N136").Value = -24718.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1153.25
$ws.Range("I20").Value = 1003.6667
$ws.Range("K20").Value = 1003.6667
$ws.Range("M20").Value = -756.6667
$ws.Range("H64").Value = 980.6667
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 980.6667
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H80").Value = 978.8570999999999
$ws.Range("I80").Value = 558.6667
$ws.Range("K80").Value = 558.6667
$ws.Range("M80").Value = 439.3333
$ws.Range("H83").Value = 978.8570999999999
$ws.Range("I83").Value = 558.6667
$ws.Range("K83").Value = 2793.3335
$ws.Range("M83").Value = 2198.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6767.9653
$ws.Range("I31").Value = 3426
$ws.Range("K31").Value = 3426
$ws.Range("M31").Value = -3131
$ws.Range("H34").Value = 6767.9653
$ws.Range("I34").Value = 3426
$ws.Range("K34").Value = 3426
$ws.Range("M34").Value = -3224
$ws.Range("H58").Value = 2229.818
$ws.Range("I58").Value = 1374.1666
$ws.Range("K58").Value = 1374.1666
$ws.Range("M58").Value = -1171.1666
$ws.Range("H69").Value = 23419.6
$ws.Range("I69").Value = 13965.667
$ws.Range("K69").Value = 13965.667
$ws.Range("M69").Value = -13216.667
$ws.Range("H72").Value = 23419.6
$ws.Range("I72").Value = 13965.667
$ws.Range("K72").Value = 41897.001
$ws.Range("M72").Value = -38153.001
$ws.Range("H105").Value = 1570.375
$ws.Range("I105").Value = 929.875
$ws.Range("K105").Value = 929.875
$ws.Range("M105").Value = 817.125
$ws.Range("H136").Value = 2229.818
$ws.Range("I136").Value = 1374.1666
$ws.Range("K136").Value = 4122.4998
$ws.Range("M136").Value = -1572.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 138.27272
$ws.Range("I6").Value = 53.75
$ws.Range("J6").Value = 363.66666
$ws.Range("K6").Value = 161.25
$ws.Range("L6").Value = 1090.99998
$ws.Range("M6").Value = -48.25
$ws.Range("N6").Value = -1316.99998
$ws.Range("H12").Value = 146.53847
$ws.Range("I12").Value = 29
$ws.Range("K12").Value = 87
$ws.Range("M12").Value = 86
$ws.Range("H40").Value = 275.5
$ws.Range("I40").Value = 5.5
$ws.Range("J40").Value = 410.5
$ws.Range("K40").Value = 22
$ws.Range("L40").Value = 1642
$ws.Range("M40").Value = 47
$ws.Range("N40").Value = -1780
$ws.Range("H81").Value = 1638.5
$ws.Range("J81").Value = 2166.3333
$ws.Range("L81").Value = 6498.999899999999
$ws.Range("N81").Value = -8744.999899999999
$ws.Range("H84").Value = 1638.5
$ws.Range("J84").Value = 2166.3333
$ws.Range("L84").Value = 19496.9997
$ws.Range("N84").Value = -30728.9997
$ws.Range("H104").Value = 9158.25
$ws.Range("J104").Value = 9989.9
$ws.Range("L104").Value = 29969.7
$ws.Range("N104").Value = -35211.7
$ws.Range("H132").Value = 1642.1428
$ws.Range("I132").Value = 1699
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 15291
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -12761
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 91.8125
$ws.Range("I2").Value = 47.5
$ws.Range("J2").Value = 165.66667
$ws.Range("K2").Value = 47.5
$ws.Range("L2").Value = 165.66667
$ws.Range("M2").Value = 65.5
$ws.Range("N2").Value = -391.66667
$ws.Range("H31").Value = 251.8
$ws.Range("I31").Value = 251.8
$ws.Range("K31").Value = 251.8
$ws.Range("M31").Value = 40.19999999999999
$ws.Range("H37").Value = 251.8
$ws.Range("I37").Value = 251.8
$ws.Range("K37").Value = 251.8
$ws.Range("M37").Value = 25.19999999999999
$ws.Range("H122").Value = 2797
$ws.Range("I122").Value = 2595
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 7785
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -5335
$ws.Range("N122").Value = -13897
$ws.Range("H126").Value = 2598.8
$ws.Range("I126").Value = 2598.8
$ws.Range("K126").Value = 7796.400000000001
$ws.Range("M126").Value = -5326.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1001
$ws.Range("I16").Value = 1001
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1001
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -831
$ws.Range("N16").ClearContents()
$ws.Range("H42").Value = 29999
$ws.Range("I42").Value = 29999
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 29999
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -29436
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 29999
$ws.Range("I49").Value = 29999
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 29999
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -29852
$ws.Range("N49").ClearContents()
$ws.Range("H125").Value = 190357.5
$ws.Range("J125").Value = 190357.5
$ws.Range("L125").Value = 190357.5
$ws.Range("N125").Value = -200197.5
$ws.Range("H136").Value = 3187.75
$ws.Range("I136").Value = 2357.4285
$ws.Range("K136").Value = 7072.2855
$ws.Range("M136").Value = -4522.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11305.9375
$ws.Range("I62").Value = 8974.5
$ws.Range("J62").Value = 12083.083
$ws.Range("K62").Value = 8974.5
$ws.Range("L62").Value = 12083.083
$ws.Range("M62").Value = -8350.5
$ws.Range("N62").Value = -13331.083
$ws.Range("H65").Value = 11305.9375
$ws.Range("I65").Value = 8974.5
$ws.Range("J65").Value = 12083.083
$ws.Range("K65").Value = 44872.5
$ws.Range("L65").Value = 60415.415
$ws.Range("M65").Value = -41752.5
$ws.Range("N65").Value = -66655.41500000001
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 6113.846
$ws.Range("I126").Value = 4013.3333
$ws.Range("K126").Value = 12039.9999
$ws.Range("M126").Value = -9569.999899999999
$ws.Range("H136").Value = 3199.4333
$ws.Range("I136").Value = 2271
$ws.Range("J136").Value = 5056.3
$ws.Range("K136").Value = 6813
$ws.Range("M136").Value = -4263
$ws.Range("N136").Value = -20268.9
